$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "doSearch"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "doSearch"

# Fill in the vehicle configuration data.
# Shared-string insertion order matters for an exact OOXML match, so write
# the model names (and values) before the header row text.
$ws.Range("A2").Value = "BMW M4"
$ws.Range("A3").Value = "BMW M8"
$ws.Range("A4").Value = "BMW X2"
$ws.Range("A5").Value = "BMW X5"
$ws.Range("A6").Value = "BMW i4"

$ws.Range("B1").Value = "Year"
$ws.Range("A1").Value = "Model"

$ws.Range("B2").Value = 2023
$ws.Range("B3").Value = 2023
$ws.Range("B4").Value = 2023
$ws.Range("B5").Value = 2023
$ws.Range("B6").Value = 2022

# Match the column-A width set on the source sheet
$ws.Columns.Item(1).ColumnWidth = 10

# Make doSearch the active sheet/tab
$ws.Activate()
